$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = -7.517499999999993
$ws.Range("C12").Value = -11.3786
$ws.Range("D23").Value = -8.029200000000005
$ws.Range("D28").Value = -8.128699999999998
$ws.Range("C32").Value = -13.1992
$ws.Range("D32").Value = -8.178599999999998
$ws.Range("D34").Value = -7.869999999999999
$ws.Range("C36").Value = -12.4061
$ws.Range("C38").Value = -12.49539999999999
$ws.Range("D42").Value = -8.932699999999992
$ws.Range("C46").Value = -14.63669999999999
$ws.Range("C54").Value = -12.5197
$ws.Range("D54").Value = -8.086700000000006
$ws.Range("C55").Value = -13.96159999999999
$ws.Range("C67").Value = -12.357
$ws.Range("C69").Value = -11.9468
$ws.Range("C72").Value = -11.9806
$ws.Range("C91").Value = -12.23280000000001
$ws.Range("D97").Value = -8.501199999999994
$ws.Range("C99").Value = -12.0548
$ws.Range("D99").Value = -8.432999999999993
$ws.Range("D101").Value = -7.8097
$ws.Range("C104").Value = -12.65660000000001
